# Re-insert the leading "+" marker on four footnotes (the footnote text
# runs that cite the Peking ("pe cin") edition). Each of these footnote
# texts gets a "+" prepended to its wording.
#
# Note: Find/Replace on a footnote's Range does not reliably commit text
# changes in this runtime, so we rebuild each target footnote's text
# (after the footnote-reference mark) via the Range.Text property, which
# is the form of mutation that persists for footnote content.

$d = $word.ActiveDocument

$updates = @{
    1  = "ཨུ་པ་དྷེ་ཤ། པེ་ཅིན།"
    5  = "ཨ། པེ་ཅིན།"
    7  = "ཅིང། པེ་ཅིན།"
    19 = "།།་མངྒཱ་ལཾ། པེ་ཅིན།"
}

foreach ($index in $updates.Keys) {
    $fn = $d.Footnotes.Item($index)
    $newText = "+" + $updates[$index]
    $fn.Range.Text = $newText
}
